$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.843.24"
$ws.Range("E2").Value = "  -1.39%  "

$ws.Range("D3").Value = "2.397.85"
$ws.Range("E3").Value = "  -2.44%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "552.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.95%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.86%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.504"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("E9").Value = "  +4.50%  "

$ws.Range("E10").Value = "  -1.36%  "

$ws.Range("E11").Value = "  -2.48%  "

$ws.Range("E12").Value = "  -2.55%  "

$ws.Range("D13").Value = "67.734.64"
$ws.Range("E13").Value = "  -1.41%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000169"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.27%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.51%  "

$ws.Range("E16").Value = "  -4.80%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "329.54"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.13%  "

$ws.Range("E18").Value = "  -3.83%  "

$ws.Range("E19").Value = "  -1.02%  "

$ws.Range("E20").Value = "  -0.09%  "

$ws.Range("E21").Value = "  -5.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.12%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.62"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.68%  "

$ws.Range("E25").Value = "  -3.31%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.37%  "

$ws.Range("E27").Value = "  -0.06%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "418.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.67%  "

$ws.Range("E29").Value = "  -1.28%  "

$ws.Range("E30").Value = "  -1.96%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "158.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.30%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.39%  "

$ws.Range("E33").Value = "  -0.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.61"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.14%  "

$ws.Range("E35").Value = "  -4.26%  "

$ws.Range("E36").Value = "  -2.99%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.22"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.27%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.05"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.43%  "

$ws.Range("E40").Value = "  -2.17%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "128.23"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.56%  "

$ws.Range("E42").Value = "  -7.59%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0705"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.95%  "

$ws.Range("E44").Value = "  -1.93%  "

$ws.Range("E45").Value = "  -1.59%  "

$ws.Range("E47").Value = "  -1.09%  "

$ws.Range("E48").Value = "  -6.42%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.18%  "

$ws.Range("D50").Value = "0.0₆0202"
$ws.Range("E50").Value = "  -6.80%  "

$ws.Range("E51").Value = "  -0.87%  "
